# Auto-generated update of market/price data in each profession sheet's
# Leve profit table (columns H-N). Mirrors a scheduled market-data refresh:
# plain numeric values only, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 315.3846
$ws.Range("I18").Value = 315.3846
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 315.3846
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -31.38459999999998
$ws.Range("N18").ClearContents()
$ws.Range("H74").Value = 4572.4
$ws.Range("I74").Value = 4183.3335
$ws.Range("J74").Value = 4831.778
$ws.Range("K74").Value = 4183.3335
$ws.Range("L74").Value = 4831.778
$ws.Range("M74").Value = -3247.3335
$ws.Range("N74").Value = -6703.778
$ws.Range("H77").Value = 4572.4
$ws.Range("I77").Value = 4183.3335
$ws.Range("J77").Value = 4831.778
$ws.Range("K77").Value = 20916.6675
$ws.Range("L77").Value = 24158.89
$ws.Range("M77").Value = -16236.6675
$ws.Range("N77").Value = -33518.89
$ws.Range("H86").Value = 39142.15
$ws.Range("I86").Value = 101888.6
$ws.Range("J86").Value = 2232.4707
$ws.Range("K86").Value = 101888.6
$ws.Range("L86").Value = 2232.4707
$ws.Range("M86").Value = -100765.6
$ws.Range("N86").Value = -4478.4707
$ws.Range("H89").Value = 39142.15
$ws.Range("I89").Value = 101888.6
$ws.Range("J89").Value = 2232.4707
$ws.Range("K89").Value = 509443
$ws.Range("L89").Value = 11162.3535
$ws.Range("M89").Value = -503827
$ws.Range("N89").Value = -22394.3535
$ws.Range("H137").Value = 2407.476
$ws.Range("I137").Value = 1395.8
$ws.Range("J137").Value = 3672.0715
$ws.Range("K137").Value = 4187.4
$ws.Range("L137").Value = 11016.2145
$ws.Range("M137").Value = -1637.4
$ws.Range("N137").Value = -16116.2145

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 720.35596
$ws.Range("I2").Value = 607.0217
$ws.Range("J2").Value = 1121.3846
$ws.Range("K2").Value = 607.0217
$ws.Range("L2").Value = 1121.3846
$ws.Range("M2").Value = -494.0217
$ws.Range("N2").Value = -1347.3846
$ws.Range("H21").Value = 1003.75
$ws.Range("I21").Value = 1003.75
$ws.Range("K21").Value = 1003.75
$ws.Range("M21").Value = -629.75
$ws.Range("H32").Value = 6333.122
$ws.Range("I32").Value = 4925.5913
$ws.Range("J32").Value = 11592.842
$ws.Range("K32").Value = 4925.5913
$ws.Range("L32").Value = 11592.842
$ws.Range("M32").Value = -4638.5913
$ws.Range("N32").Value = -12166.842
$ws.Range("H116").Value = 720.35596
$ws.Range("I116").Value = 607.0217
$ws.Range("J116").Value = 1121.3846
$ws.Range("K116").Value = 607.0217
$ws.Range("L116").Value = 1121.3846
$ws.Range("M116").Value = 1686.9783
$ws.Range("N116").Value = -5709.3846
$ws.Range("H122").Value = 2520.6365
$ws.Range("I122").Value = 2727.2307
$ws.Range("J122").Value = 2222.2222
$ws.Range("K122").Value = 8181.6921
$ws.Range("L122").Value = 6666.6666
$ws.Range("M122").Value = -5731.6921
$ws.Range("N122").Value = -11566.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 720.35596
$ws.Range("I3").Value = 607.0217
$ws.Range("J3").Value = 1121.3846
$ws.Range("K3").Value = 607.0217
$ws.Range("L3").Value = 1121.3846
$ws.Range("M3").Value = -493.0217
$ws.Range("N3").Value = -1349.3846
$ws.Range("H105").Value = 5630.1816
$ws.Range("I105").Value = 5616.5557
$ws.Range("J105").Value = 5646.533
$ws.Range("K105").Value = 5616.5557
$ws.Range("L105").Value = 5646.533
$ws.Range("M105").Value = -3869.5557
$ws.Range("N105").Value = -9140.532999999999
$ws.Range("H126").Value = 71664
$ws.Range("J126").Value = 71664
$ws.Range("L126").Value = 71664
$ws.Range("N126").Value = -81544
$ws.Range("H141").Value = 41256.668
$ws.Range("J141").Value = 41256.668
$ws.Range("L141").Value = 41256.668
$ws.Range("N141").Value = -51616.668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1718.25
$ws.Range("I31").Value = 1316.2632
$ws.Range("K31").Value = 1316.2632
$ws.Range("M31").Value = -1021.2632
$ws.Range("H34").Value = 1718.25
$ws.Range("I34").Value = 1316.2632
$ws.Range("K34").Value = 1316.2632
$ws.Range("M34").Value = -1114.2632
$ws.Range("H58").Value = 1857889.9
$ws.Range("I58").Value = 3136513
$ws.Range("J58").Value = 3886.1
$ws.Range("K58").Value = 3136513
$ws.Range("L58").Value = 3886.1
$ws.Range("M58").Value = -3136310
$ws.Range("N58").Value = -4292.1
$ws.Range("H99").Value = 1408.174
$ws.Range("I99").Value = 1391.3125
$ws.Range("K99").Value = 1391.3125
$ws.Range("M99").Value = 106.6875
$ws.Range("H122").Value = 12143.125
$ws.Range("I122").Value = 4004.125
$ws.Range("K122").Value = 12012.375
$ws.Range("M122").Value = -9562.375
$ws.Range("H126").Value = 1408.174
$ws.Range("I126").Value = 1391.3125
$ws.Range("K126").Value = 4173.9375
$ws.Range("M126").Value = -1703.9375
$ws.Range("H136").Value = 1857889.9
$ws.Range("I136").Value = 3136513
$ws.Range("J136").Value = 3886.1
$ws.Range("K136").Value = 9409539
$ws.Range("L136").Value = 11658.3
$ws.Range("M136").Value = -9406989
$ws.Range("N136").Value = -16758.3

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 23447.809
$ws.Range("I131").Value = 651.3158
$ws.Range("J131").Value = 42279.695
$ws.Range("K131").Value = 1953.9474
$ws.Range("L131").Value = 126839.085
$ws.Range("M131").Value = 3086.0526
$ws.Range("N131").Value = -136919.085
$ws.Range("H132").Value = 2249.0833
$ws.Range("I132").Value = 2483.1667
$ws.Range("J132").Value = 2015
$ws.Range("K132").Value = 22348.5003
$ws.Range("L132").Value = 18135
$ws.Range("M132").Value = -19818.5003
$ws.Range("N132").Value = -23195

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 22037
$ws.Range("J45").Value = 22037
$ws.Range("L45").Value = 22037
$ws.Range("N45").Value = -23155
$ws.Range("H126").Value = 3400.5881
$ws.Range("I126").Value = 2346.3635
$ws.Range("J126").Value = 5333.3335
$ws.Range("K126").Value = 7039.0905
$ws.Range("L126").Value = 16000.0005
$ws.Range("M126").Value = -4569.0905
$ws.Range("N126").Value = -20940.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5104.773
$ws.Range("I7").Value = 3317.0833
$ws.Range("J7").Value = 7250
$ws.Range("K7").Value = 3317.0833
$ws.Range("L7").Value = 7250
$ws.Range("M7").Value = -3205.0833
$ws.Range("N7").Value = -7474
$ws.Range("H40").Value = 4363.5
$ws.Range("I40").Value = 3999.8125
$ws.Range("J40").Value = 5333.3335
$ws.Range("K40").Value = 3999.8125
$ws.Range("L40").Value = 5333.3335
$ws.Range("M40").Value = -3863.8125
$ws.Range("N40").Value = -5605.3335
$ws.Range("H61").Value = 18755.666
$ws.Range("I61").Value = 26866.416
$ws.Range("J61").Value = 2534.1667
$ws.Range("K61").Value = 26866.416
$ws.Range("L61").Value = 2534.1667
$ws.Range("M61").Value = -26664.416
$ws.Range("N61").Value = -2938.1667
$ws.Range("H113").Value = 18755.666
$ws.Range("I113").Value = 26866.416
$ws.Range("J113").Value = 2534.1667
$ws.Range("K113").Value = 26866.416
$ws.Range("L113").Value = 2534.1667
$ws.Range("M113").Value = -24696.416
$ws.Range("N113").Value = -6874.1667
$ws.Range("H122").Value = 4865.7812
$ws.Range("I122").Value = 4395.4546
$ws.Range("J122").Value = 5900.5
$ws.Range("K122").Value = 13186.3638
$ws.Range("L122").Value = 17701.5
$ws.Range("M122").Value = -10736.3638
$ws.Range("N122").Value = -22601.5
$ws.Range("H126").Value = 5104.773
$ws.Range("I126").Value = 3317.0833
$ws.Range("J126").Value = 7250
$ws.Range("K126").Value = 9951.249899999999
$ws.Range("L126").Value = 21750
$ws.Range("M126").Value = -7481.249899999999
$ws.Range("N126").Value = -26690
$ws.Range("H136").Value = 4880.655
$ws.Range("I136").Value = 3629.7334
$ws.Range("K136").Value = 10889.2002
$ws.Range("M136").Value = -8339.200199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1222.091
$ws.Range("I107").Value = 356.81818
$ws.Range("J107").Value = 2087.3635
$ws.Range("K107").Value = 1070.45454
$ws.Range("L107").Value = 6262.0905
$ws.Range("M107").Value = 849.54546
$ws.Range("N107").Value = -10102.0905
$ws.Range("H113").Value = 1124.4584
$ws.Range("I113").Value = 512.06665
$ws.Range("J113").Value = 2145.111
$ws.Range("K113").Value = 1536.19995
$ws.Range("L113").Value = 6435.333
$ws.Range("M113").Value = 633.8000500000001
$ws.Range("N113").Value = -10775.333
$ws.Range("H136").Value = 3917.7014
$ws.Range("I136").Value = 1727.1904
$ws.Range("J136").Value = 7597.76
$ws.Range("K136").Value = 5181.5712
$ws.Range("L136").Value = 22793.28
$ws.Range("M136").Value = -2631.5712
$ws.Range("N136").Value = -27893.28
